# Commit: "fix lỗi trong report cơ sở. Thêm cột ghi chú trong báo cáo về chi tiêu"
# This populates the (previously empty) "Đơn sale chính" sheet with a header row,
# one order row and a totals row, and updates the pre-computed figures on the
# "Lương" sheet to reflect the new order.

$wb = $excel.ActiveWorkbook

# Force a value to land as a literal (inline) string, never letting Excel's
# "looks like a date/number" auto-conversion kick in: write it as a formula
# returning the literal text, then flatten formula -> value via copy/paste.
function SetText($sheet, $addr, $text) {
    $rng = $sheet.Range($addr)
    $rng.Formula = $text
    $rng.Copy()
    $rng.PasteSpecial(-4163)
}

# ---------------------------------------------------------------------------
# Sheet 1: "Đơn sale chính"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Row 1 - headers
SetText $ws1 "A1" '="Tiền tố"'
SetText $ws1 "B1" '="Mã dịch vụ"'
SetText $ws1 "C1" '="Ngày thực hiện"'
SetText $ws1 "D1" '="Cơ sở"'
SetText $ws1 "E1" '="Khách hàng"'
SetText $ws1 "F1" '="Nguồn khách"'
SetText $ws1 "G1" '="Tên dịch vụ"'
SetText $ws1 "H1" '="Đơn giá gốc"'
SetText $ws1 "I1" '="Sale phụ"'
SetText $ws1 "J1" '="Upsale"'
SetText $ws1 "K1" '="Đơn giá"'
SetText $ws1 "L1" '="Đã thanh toán"'
SetText $ws1 "M1" '="Tỉ lệ chiết khấu sale chính"'
SetText $ws1 "N1" '="Chiết khấu sale chính"'

# Row 2 - the order
SetText $ws1 "A2" '="HD-LUXURY"'
$ws1.Range("B2").Value = 632
SetText $ws1 "C2" '="08-07-2024"'
SetText $ws1 "D2" '="CẦN THƠ"'
SetText $ws1 "E2" '="Nguyễn Thị Thắm"'
SetText $ws1 "F2" '="Cá nhân"'
SetText $ws1 "G2" '="Tiêm botox"'
$ws1.Range("H2").Value = 2000000
$ws1.Range("I2").Value = 0
$ws1.Range("J2").Value = 0
$ws1.Range("K2").Value = 2000000
$ws1.Range("L2").Value = 2000000
$ws1.Range("M2").Value = 0.1
$ws1.Range("N2").Value = 200000

# Row 3 - totals
SetText $ws1 "A3" '="Tổng"'
$ws1.Range("B3").Value = 1
$ws1.Range("H3").Value = 2000000
$ws1.Range("J3").Value = 0
$ws1.Range("K3").Value = 2000000
$ws1.Range("L3").Value = 2000000
$ws1.Range("M3").Value = 0
$ws1.Range("N3").Value = 200000

# ---------------------------------------------------------------------------
# Sheet 3: "Lương" - figures recomputed after the new order above
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("B2").Value = 8
$ws3.Range("B3").Value = 280000
$ws3.Range("B4").Value = 4285714.285714285
$ws3.Range("B5").Value = 200000
$ws3.Range("B15").Value = 2857142.857142857
$ws3.Range("B26").Value = 4285714.285714285
$ws3.Range("B35").Value = 5765714.285714285
$ws3.Range("B36").Value = 2857142.857142857
$ws3.Range("B37").Value = 4285714.285714285
$ws3.Range("B38").Value = 12908571.42857143
